$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '41.755.14'
Set-TextValue $ws.Range("E2") '  +1.44%  '
Set-TextValue $ws.Range("D3") '2.262.08'
Set-TextValue $ws.Range("E3") '  +0.67%  '
Set-TextValue $ws.Range("E4") '  -0.02%  '
Set-TextValue $ws.Range("D5") '303.43'
Set-TextValue $ws.Range("E5") '  +0.56%  '
Set-TextValue $ws.Range("D6") '92.08'
Set-TextValue $ws.Range("E6") '  +1.49%  '
Set-TextValue $ws.Range("D7") '0.532'
Set-TextValue $ws.Range("E7") '  +2.31%  '
Set-TextValue $ws.Range("E8") '  -0.05%  '
Set-TextValue $ws.Range("E9") '  +0.42%  '
Set-TextValue $ws.Range("D10") '32.44'
Set-TextValue $ws.Range("E10") '  +2.48%  '
Set-TextValue $ws.Range("D11") '53.46'
Set-TextValue $ws.Range("E11") '  -0.41%  '
Set-TextValue $ws.Range("E12") '  +0.64%  '
Set-TextValue $ws.Range("E13") '  -0.46%  '
Set-TextValue $ws.Range("E14") '  +1.71%  '
Set-TextValue $ws.Range("D15") '2.616.99'
Set-TextValue $ws.Range("E15") '  +0.91%  '
Set-TextValue $ws.Range("D16") '14.26'
Set-TextValue $ws.Range("E16") '  +1.45%  '
Set-TextValue $ws.Range("D17") '2.282.10'
Set-TextValue $ws.Range("E17") '  +4.07%  '
Set-TextValue $ws.Range("D18") '0.773'
Set-TextValue $ws.Range("E18") '  +3.26%  '
Set-TextValue $ws.Range("D19") '41.655.51'
Set-TextValue $ws.Range("E19") '  +1.37%  '
Set-TextValue $ws.Range("D20") '12.40'
Set-TextValue $ws.Range("E20") '  +5.09%  '
Set-TextValue $ws.Range("E21") '  +0.63%  '
Set-TextValue $ws.Range("E22") '  +1.68%  '
Set-TextValue $ws.Range("D23") '67.10'
Set-TextValue $ws.Range("E23") '  +0.60%  '
Set-TextValue $ws.Range("D24") '239.73'
Set-TextValue $ws.Range("E24") '  +0.11%  '
Set-TextValue $ws.Range("D25") '2.58'
Set-TextValue $ws.Range("E25") '  +1.27%  '
Set-TextValue $ws.Range("E26") '  +0.04%  '
Set-TextValue $ws.Range("E27") '  +4.07%  '
Set-TextValue $ws.Range("D28") '23.93'
Set-TextValue $ws.Range("E28") '  +0.96%  '
Set-TextValue $ws.Range("D29") '9.52'
Set-TextValue $ws.Range("E29") '  -0.01%  '
Set-TextValue $ws.Range("E30") '  -4.45%  '
Set-TextValue $ws.Range("D31") '35.05'
Set-TextValue $ws.Range("E31") '  +6.19%  '
Set-TextValue $ws.Range("D32") '160.50'
Set-TextValue $ws.Range("E32") '  +1.17%  '
Set-TextValue $ws.Range("D33") '5.26'
Set-TextValue $ws.Range("E33") '  +2.15%  '
Set-TextValue $ws.Range("E34") '  -0.08%  '
Set-TextValue $ws.Range("E35") '  +1.82%  '
Set-TextValue $ws.Range("E36") '  -0.12%  '
Set-TextValue $ws.Range("D37") '16.96'
Set-TextValue $ws.Range("E37") '  +3.84%  '
Set-TextValue $ws.Range("E38") '  +0.23%  '
Set-TextValue $ws.Range("E39") '  +1.51%  '
Set-TextValue $ws.Range("D40") '0.105'
Set-TextValue $ws.Range("E40") '  +0.84%  '
Set-TextValue $ws.Range("E41") '  +0.94%  '
Set-TextValue $ws.Range("E42") '  +0.19%  '
Set-TextValue $ws.Range("D43") '2.014.03'
Set-TextValue $ws.Range("E43") '  -2.87%  '
Set-TextValue $ws.Range("D44") '19.27'
Set-TextValue $ws.Range("E44") '  -4.23%  '
Set-TextValue $ws.Range("E45") '  +1.03%  '
Set-TextValue $ws.Range("D46") '10.29'
Set-TextValue $ws.Range("E46") '  +0.79%  '
Set-TextValue $ws.Range("E47") '  +5.39%  '
Set-TextValue $ws.Range("E48") '  -2.21%  '
Set-TextValue $ws.Range("E49") '  +1.13%  '
Set-TextValue $ws.Range("E50") '  +1.17%  '

# Row 51: MultiversX -> BitcoinSV (full row replacement)
Set-TextValue $ws.Range("B51") 'BitcoinSV'
Set-TextValue $ws.Range("C51") 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range("D51") '72.09'
Set-TextValue $ws.Range("E51") '  +2.73%  '
